$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 136.2
$ws.Range("I8").Value = 136.2
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 408.6
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -269.6
$ws.Range("N8").Value = $null

$ws.Range("H135").Value = 2348.9
$ws.Range("I135").Value = 1498.7778
$ws.Range("J135").Value = 10000
$ws.Range("K135").Value = 13489.0002
$ws.Range("L135").Value = 90000
$ws.Range("M135").Value = -10954.0002
$ws.Range("N135").Value = -95070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 39250
$ws.Range("J24").Value = 39250
$ws.Range("L24").Value = 39250
$ws.Range("N24").Value = -39998

$ws.Range("H74").Value = 3723.8164
$ws.Range("I74").Value = 897.25
$ws.Range("J74").Value = 9044.412
$ws.Range("K74").Value = 897.25
$ws.Range("L74").Value = 9044.412
$ws.Range("M74").Value = -23.25
$ws.Range("N74").Value = -10792.412

$ws.Range("H77").Value = 3723.8164
$ws.Range("I77").Value = 897.25
$ws.Range("J77").Value = 9044.412
$ws.Range("K77").Value = 4486.25
$ws.Range("L77").Value = 45222.06
$ws.Range("M77").Value = -118.25
$ws.Range("N77").Value = -53958.06

$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492

$ws.Range("H100").Value = 39250
$ws.Range("J100").Value = 39250
$ws.Range("L100").Value = 39250
$ws.Range("N100").Value = -41414

$ws.Range("H132").Value = 1610.4375
$ws.Range("I132").Value = 1647.1957
$ws.Range("K132").Value = 4941.5871
$ws.Range("M132").Value = -2411.5871

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 29600
$ws.Range("J92").Value = 29600
$ws.Range("L92").Value = 29600
$ws.Range("N92").Value = -34592

$ws.Range("H96").Value = 32924.5
$ws.Range("J96").Value = 32924.5
$ws.Range("L96").Value = 32924.5
$ws.Range("N96").Value = -38416.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1530.3334
$ws.Range("I5").Value = 1018.8571
$ws.Range("J5").Value = 1653.7931
$ws.Range("K5").Value = 3056.5713
$ws.Range("L5").Value = 4961.379300000001
$ws.Range("M5").Value = -2944.5713
$ws.Range("N5").Value = -5185.379300000001

$ws.Range("H64").Value = 3699.5293
$ws.Range("I64").Value = 1397.6
$ws.Range("J64").Value = 4658.6665
$ws.Range("K64").Value = 4192.799999999999
$ws.Range("L64").Value = 13975.9995
$ws.Range("M64").Value = -3922.799999999999
$ws.Range("N64").Value = -14515.9995

$ws.Range("H67").Value = 3699.5293
$ws.Range("I67").Value = 1397.6
$ws.Range("J67").Value = 4658.6665
$ws.Range("K67").Value = 4192.799999999999
$ws.Range("L67").Value = 13975.9995
$ws.Range("M67").Value = -3256.799999999999
$ws.Range("N67").Value = -15847.9995

$ws.Range("H135").Value = 1530.3334
$ws.Range("I135").Value = 1018.8571
$ws.Range("J135").Value = 1653.7931
$ws.Range("K135").Value = 9169.713899999999
$ws.Range("L135").Value = 14884.1379
$ws.Range("M135").Value = -6634.713899999999
$ws.Range("N135").Value = -19954.1379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 29116.834
$ws.Range("J39").Value = 29116.834
$ws.Range("L39").Value = 29116.834
$ws.Range("N39").Value = -30180.834

$ws.Range("H70").Value = 15633.333
$ws.Range("I70").Value = 19950
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 19950
$ws.Range("L70").Value = 7000
$ws.Range("M70").Value = -19680
$ws.Range("N70").Value = -7540

$ws.Range("H73").Value = 15633.333
$ws.Range("I73").Value = 19950
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 19950
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = -19014
$ws.Range("N73").Value = -8872

$ws.Range("H102").Value = 3270.5
$ws.Range("I102").Value = 4240.125
$ws.Range("J102").Value = 2624.0833
$ws.Range("K102").Value = 4240.125
$ws.Range("L102").Value = 2624.0833
$ws.Range("M102").Value = -2618.125
$ws.Range("N102").Value = -5868.0833

$ws.Range("H113").Value = 1229.375
$ws.Range("I113").Value = 1105.5
$ws.Range("J113").Value = 1353.25
$ws.Range("K113").Value = 1105.5
$ws.Range("L113").Value = 1353.25
$ws.Range("M113").Value = 1064.5
$ws.Range("N113").Value = -5693.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1539
$ws.Range("I22").Value = 1700
$ws.Range("J22").Value = 1498.75
$ws.Range("K22").Value = 1700
$ws.Range("L22").Value = 1498.75
$ws.Range("M22").Value = -1405
$ws.Range("N22").Value = -2088.75

$ws.Range("H27").Value = 1539
$ws.Range("I27").Value = 1700
$ws.Range("J27").Value = 1498.75
$ws.Range("K27").Value = 1700
$ws.Range("L27").Value = 1498.75
$ws.Range("M27").Value = -1593
$ws.Range("N27").Value = -1712.75

$ws.Range("H81").Value = 32681
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 32681
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 32681
$ws.Range("M81").Value = $null
$ws.Range("N81").Value = -34677

$ws.Range("H84").Value = 32681
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 32681
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 98043
$ws.Range("M84").Value = $null
$ws.Range("N84").Value = -108027

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J33").Value = 50000
$ws.Range("L33").Value = 50000
$ws.Range("N33").Value = -50500

$ws.Range("J36").Value = 50000
$ws.Range("L36").Value = 50000
$ws.Range("N36").Value = -50500

$ws.Range("H97").Value = 30000
$ws.Range("J97").Value = 30000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -31982

$ws.Range("H132").Value = 3990.6
$ws.Range("I132").Value = 4242.3335
$ws.Range("J132").Value = 1725
$ws.Range("K132").Value = 12727.0005
$ws.Range("L132").Value = 5175
$ws.Range("M132").Value = -10197.0005
$ws.Range("N132").Value = -10235
